# Applies the weekly re-shuffle of the Hortaliza / Sandia price rows.
# Each target row below is overwritten with the values that originally lived
# on a different row of the same sheet (date, quality grade, volume, min/max/
# weighted-avg price, origin, and $/Kg), matching the "Fruta / hortaliza,
# semanal" commit. Rows 5, 6, 7, 27 and 28 are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <- original row 18
$ws.Range("D2").Value = 44214
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 1200
$ws.Range("K2").Value = 400
$ws.Range("L2").Value = 450
$ws.Range("M2").Value = 425
$ws.Range("O2").Value = "Perú"
$ws.Range("P2").Value = 425

# Row 3 <- original row 19
$ws.Range("D3").Value = 44172
$ws.Range("I3").Value = "Segunda"
$ws.Range("J3").Value = 1600
$ws.Range("K3").Value = 400
$ws.Range("L3").Value = 420
$ws.Range("M3").Value = 410
$ws.Range("O3").Value = "Perú"
$ws.Range("P3").Value = 410

# Row 4 <- original row 9
$ws.Range("D4").Value = 44575
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 1200
$ws.Range("K4").Value = 380
$ws.Range("L4").Value = 400
$ws.Range("M4").Value = 390
$ws.Range("O4").Value = "Perú"
$ws.Range("P4").Value = 390

# Row 8 <- original row 20
$ws.Range("D8").Value = 44201
$ws.Range("I8").Value = "Segunda"
$ws.Range("J8").Value = 1800
$ws.Range("K8").Value = 250
$ws.Range("L8").Value = 270
$ws.Range("M8").Value = 260
$ws.Range("O8").Value = "Perú"
$ws.Range("P8").Value = 260

# Row 9 <- original row 25
$ws.Range("D9").Value = 44229
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 1200
$ws.Range("K9").Value = 230
$ws.Range("L9").Value = 250
$ws.Range("M9").Value = 240
$ws.Range("O9").Value = "Perú"
$ws.Range("P9").Value = 240

# Row 10 <- original row 26
$ws.Range("D10").Value = 44301
$ws.Range("I10").Value = "Segunda"
$ws.Range("J10").Value = 900
$ws.Range("K10").Value = 280
$ws.Range("L10").Value = 300
$ws.Range("M10").Value = 290
$ws.Range("O10").Value = "Perú"
$ws.Range("P10").Value = 290

# Row 11 <- original row 24
$ws.Range("D11").Value = 44251
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 1200
$ws.Range("K11").Value = 250
$ws.Range("L11").Value = 280
$ws.Range("M11").Value = 265
$ws.Range("O11").Value = "Perú"
$ws.Range("P11").Value = 265

# Row 12 <- original row 2
$ws.Range("D12").Value = 44224
$ws.Range("I12").Value = "Segunda"
$ws.Range("J12").Value = 1200
$ws.Range("K12").Value = 230
$ws.Range("L12").Value = 250
$ws.Range("M12").Value = 240
$ws.Range("O12").Value = "Perú"
$ws.Range("P12").Value = 240

# Row 13 <- original row 3
$ws.Range("D13").Value = 44224
$ws.Range("I13").Value = "Segunda"
$ws.Range("J13").Value = 200
$ws.Range("K13").Value = 200
$ws.Range("L13").Value = 230
$ws.Range("M13").Value = 215
$ws.Range("O13").Value = "Región de Arica y Parinacota"
$ws.Range("P13").Value = 215

# Row 14 <- original row 22
$ws.Range("D14").Value = 44168
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 1700
$ws.Range("K14").Value = 430
$ws.Range("L14").Value = 450
$ws.Range("M14").Value = 440
$ws.Range("O14").Value = "Perú"
$ws.Range("P14").Value = 440

# Row 15 <- original row 13
$ws.Range("D15").Value = 44231
$ws.Range("I15").Value = "Segunda"
$ws.Range("J15").Value = 200
$ws.Range("K15").Value = 180
$ws.Range("L15").Value = 200
$ws.Range("M15").Value = 190
$ws.Range("O15").Value = "Región de Arica y Parinacota"
$ws.Range("P15").Value = 190

# Row 16 <- original row 8
$ws.Range("D16").Value = 44547
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 1200
$ws.Range("K16").Value = 350
$ws.Range("L16").Value = 370
$ws.Range("M16").Value = 360
$ws.Range("O16").Value = "Perú"
$ws.Range("P16").Value = 360

# Row 17 <- original row 15
$ws.Range("D17").Value = 44530
$ws.Range("I17").Value = "Segunda"
$ws.Range("J17").Value = 1300
$ws.Range("K17").Value = 450
$ws.Range("L17").Value = 480
$ws.Range("M17").Value = 465
$ws.Range("O17").Value = "Perú"
$ws.Range("P17").Value = 465

# Row 18 <- original row 14
$ws.Range("D18").Value = 44580
$ws.Range("I18").Value = "Primera"
$ws.Range("J18").Value = 1200
$ws.Range("K18").Value = 380
$ws.Range("L18").Value = 400
$ws.Range("M18").Value = 390
$ws.Range("O18").Value = "Región Metropolitana"
$ws.Range("P18").Value = 390

# Row 19 <- original row 11
$ws.Range("D19").Value = 44217
$ws.Range("I19").Value = "Segunda"
$ws.Range("J19").Value = 1600
$ws.Range("K19").Value = 300
$ws.Range("L19").Value = 350
$ws.Range("M19").Value = 325
$ws.Range("O19").Value = "Perú"
$ws.Range("P19").Value = 325

# Row 20 <- original row 12
$ws.Range("D20").Value = 44166
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value = 1700
$ws.Range("K20").Value = 500
$ws.Range("L20").Value = 530
$ws.Range("M20").Value = 515
$ws.Range("O20").Value = "Perú"
$ws.Range("P20").Value = 515

# Row 21 <- original row 16
$ws.Range("D21").Value = 44175
$ws.Range("I21").Value = "Segunda"
$ws.Range("J21").Value = 1200
$ws.Range("K21").Value = 400
$ws.Range("L21").Value = 430
$ws.Range("M21").Value = 415
$ws.Range("O21").Value = "Perú"
$ws.Range("P21").Value = 415

# Row 22 <- original row 23
$ws.Range("D22").Value = 44160
$ws.Range("I22").Value = "Segunda"
$ws.Range("J22").Value = 2000
$ws.Range("K22").Value = 500
$ws.Range("L22").Value = 550
$ws.Range("M22").Value = 525
$ws.Range("O22").Value = "Perú"
$ws.Range("P22").Value = 525

# Row 23 <- original row 10
$ws.Range("D23").Value = 44253
$ws.Range("I23").Value = "Segunda"
$ws.Range("J23").Value = 1200
$ws.Range("K23").Value = 270
$ws.Range("L23").Value = 280
$ws.Range("M23").Value = 275
$ws.Range("O23").Value = "Perú"
$ws.Range("P23").Value = 275

# Row 24 <- original row 17
$ws.Range("D24").Value = 44202
$ws.Range("I24").Value = "Segunda"
$ws.Range("J24").Value = 1300
$ws.Range("K24").Value = 230
$ws.Range("L24").Value = 250
$ws.Range("M24").Value = 240
$ws.Range("O24").Value = "Perú"
$ws.Range("P24").Value = 240

# Row 25 <- original row 21
$ws.Range("D25").Value = 44523
$ws.Range("I25").Value = "Segunda"
$ws.Range("J25").Value = 1000
$ws.Range("K25").Value = 550
$ws.Range("L25").Value = 580
$ws.Range("M25").Value = 565
$ws.Range("O25").Value = "Perú"
$ws.Range("P25").Value = 565

# Row 26 <- original row 4
$ws.Range("D26").Value = 44176
$ws.Range("I26").Value = "Primera"
$ws.Range("J26").Value = 1300
$ws.Range("K26").Value = 350
$ws.Range("L26").Value = 400
$ws.Range("M26").Value = 375
$ws.Range("O26").Value = "Perú"
$ws.Range("P26").Value = 375
